# Fixed a bug in stats2
# The data rows (A2:F20) got reshuffled - each destination row now holds
# the values that used to live in a different source row. Row 15 is
# unchanged (maps to itself). We capture the original values first, then
# write them back out in their new positions so the operation behaves
# correctly regardless of write order.
#
# Note: use Value2 (not Value) to read/write - Value misbehaves through
# this COM-interop layer and yields a descriptor string instead of the
# real contents.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: destination row -> source row (values to move into destination)
$rowMap = @{
    2  = 6
    3  = 8
    4  = 13
    5  = 3
    6  = 7
    7  = 2
    8  = 10
    9  = 11
    10 = 14
    11 = 4
    12 = 5
    13 = 9
    14 = 12
    15 = 15
    16 = 18
    17 = 20
    18 = 16
    19 = 17
    20 = 19
}

# Snapshot the original values of columns A:F for rows 2-20 before
# overwriting anything.
$original = @{}
for ($r = 2; $r -le 20; $r++) {
    $original[$r] = $ws.Range("A$r`:F$r").Value2
}

# Write each destination row using the snapshot of its source row.
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $ws.Range("A$destRow`:F$destRow").Value2 = $original[$srcRow]
}
